$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names for rows 2..26 (A column)
$names = @("model_32_8_0", "model_32_8_22", "model_32_8_21", "model_32_8_20", "model_32_8_19", "model_32_8_18", "model_32_8_17", "model_32_8_16", "model_32_8_15", "model_32_8_14", "model_32_8_13", "model_32_8_23", "model_32_8_12", "model_32_8_10", "model_32_8_9", "model_32_8_8", "model_32_8_7", "model_32_8_6", "model_32_8_5", "model_32_8_4", "model_32_8_3", "model_32_8_2", "model_32_8_1", "model_32_8_11", "model_32_8_24")

# Constant metric values (B..Q) applied identically to every row
$vals = @(0.9999949039827802, 0.9991177146462858, 0.9999966566645605, 0.9999993482079649, 0.9999983376502612, 0.000004756908385316126, 0.0008235746498410597, 0.000002922666938981576, 0.0000009453755172279311, 0.000001934021228104753, 0.00008583594873526996, 0.002181033788210565, 1.000004892176531, 0.00227388491193704, 122.5118251992398, 182.2367406177816)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $vals[$c]
    }
}
